$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 10
$ws.Range("B11").Formula = "=B10+1"
$ws.Range("C11").Value = 462
$ws.Range("D11").Value = 32
$ws.Range("E11").Value = 515
$ws.Range("F11").Value = 1009
$ws.Range("G11").Value = 354
$ws.Range("H11").Value = 526
$ws.Range("I11").Value = 0

$ws.Range("I2").Copy()
$ws.Range("C11:I11").PasteSpecial(-4122)

$ws.Range("H12").Select()
